# "moving to new server" - update the Центр label on Лист1 (sheet 1)
# and leave the selection where the user last clicked (D4).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "Центр"
$ws.Range("D4").Select() | Out-Null
